$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-05-06T13:57:20+00:00"

# --- Elements sheet: update the first data row (XCN1) ---
$elements = $wb.Worksheets.Item("Elements")

# ID / Path / Base Path columns (A, B, AF) drop the "[x]" suffix
$elements.Range("A3").Value = "ActorXDS.XCN1"
$elements.Range("B3").Value = "ActorXDS.XCN1"
$elements.Range("AF3").Value = "ActorXDS.XCN1"

# Type(s) column (K) becomes the single Identifiant URL, with a trailing newline
$elements.Range("K3").Value = "https://interop.esante.gouv.fr/ig/fhir/pdsm4dmp/StructureDefinition/Identifiant`n"

# Column K (Type(s)) now fits much shorter content - narrow its best-fit width
# (62.6484375 in the saved file, ~61.85 "characters" as reported by the COM width)
$elements.Columns.Item(11).ColumnWidth = 61.85
